$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.119.01'
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').Value = '2.621.31'
$ws.Range('E3').Value = '  +1.00%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '520.87'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.80'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.13%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.570'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.72%  '
$ws.Range('D9').Value = '2.625.65'
$ws.Range('E9').Value = '  +0.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.35'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.39%  '
$ws.Range('E11').Value = '  -0.04%  '
$ws.Range('E12').Value = '  -1.55%  '
$ws.Range('E13').Value = '  -0.48%  '
$ws.Range('D14').Value = '3.079.17'
$ws.Range('E14').Value = '  +0.95%  '
$ws.Range('D15').Value = '60.110.50'
$ws.Range('E15').Value = '  -0.78%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.23'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.00%  '
$ws.Range('E17').Value = '  -1.53%  '
$ws.Range('D18').Value = '2.626.58'
$ws.Range('E18').Value = '  +1.13%  '
$ws.Range('E19').Value = '  -2.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '341.43'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.43'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.12'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.995'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.03'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.419'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.997'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('E27').Value = '  -2.89%  '
$ws.Range('D28').Value = '0.0₃0808'
$ws.Range('E28').Value = '  -4.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.05'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.81%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.97'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.78%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.95'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '150.11'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.94'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.921'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.14'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.92%  '
$ws.Range('E38').Value = '  +2.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.56'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.54%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.43'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.63'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '288.04'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.44%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.625'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.49%  '
$ws.Range('E44').Value = '  -1.11%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('E46').Value = '  -2.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '19.39'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.63%  '
$ws.Range('E48').Value = '  +0.83%  '
$ws.Range('E49').Value = '  -2.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.66'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.05%  '
$ws.Range('D51').Value = '1.955.30'
$ws.Range('E51').Value = '  -0.15%  '
